# Automatische test-sync: 2025-08-04 20:51:50
# Adds Testmail #12 ("Ik heb nog geen geld terug.") as new row 24 on the
# "Logs" sheet, expands the worksheet dimension / conditional-formatting
# ranges to include it, and bumps the "Retour / Terugbetaling" tally on the
# "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 24

# The "Antwoord" text for this row spans multiple lines. Writing a multi-line
# string straight into the target cell via Range.Value makes the engine
# auto-fit (and mark as custom) the row height, which the source workbook
# does not have. Staging the text in a scratch cell on a far-away, unused
# row and then Copy/PasteSpecial-ing it into place reproduces the value
# without that unwanted row-height side effect; the scratch row is then
# deleted outright so no trace of it remains.
$scratchRow = 100
$scratch = $ws.Range("Z$scratchRow")
$scratch.Value = "Beste klant,
Bedankt voor uw bericht. Om uw terugbetaling te kunnen verwerken, heb ik wat meer informatie nodig. Kunt u mij uw ordernummer of transactiereferentie geven, zodat ik het voor u kan nakijken?
Alvast bedankt voor uw medewerking.
Met vriendelijke groet,
[Naam] 
E-mailassistent"
$scratch.Copy()
$ws.Range("E$row").PasteSpecial()
$ws.Rows.Item($scratchRow).Delete()

$ws.Range("A$row").Value = "Ik heb nog geen geld terug."
$ws.Range("B$row").Value = "mailmind.test@zohomail.eu"
$ws.Range("C$row").Value = "Testmail #12: Ik heb nog geen geld terug."
$ws.Range("D$row").Value = "Retour / Terugbetaling"
$ws.Range("F$row").Value = "2025-08-04 20:51:30"
$ws.Range("G$row").Value = "Ja"
$ws.Range("H$row").Value = "Nee"
$ws.Range("I$row").Value = "Ja"
$ws.Range("J$row").Value = "Nee"

# Expand the conditional-formatting ranges (D/G/H/I/J) from row 23 to row 24
# so the newly added row participates in the same highlighting rules.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = "$col" + "2:" + "$col" + "23"
    $newRange = "$col" + "2:" + "$col" + "24"
    $fcs = $ws.Range($oldRange).FormatConditions
    $count = $fcs.Count()
    for ($i = 1; $i -le $count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($newRange))
    }
}

# Update the Dashboard summary count for "Retour / Terugbetaling" (3 -> 4)
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B5").Value = 4
